$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing rows 2-11 down to 3-12.
# Excel's row insert copies formatting from the row above (the header row),
# so strip that and instead adopt the plain/date formatting used by the
# rest of the data rows.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:T2").ClearFormats()
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

# Populate the new row 2 with data (same constant columns as the rest of the
# sheet, with new date/volume/price values)
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 45111
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = "Fruta"
$ws.Cells.Item(2, 7).Value = 100104
$ws.Cells.Item(2, 8).Value = "Frutos de pepita"
$ws.Cells.Item(2, 9).Value = 100104001
$ws.Cells.Item(2, 10).Value = "Granada"
$ws.Cells.Item(2, 11).Value = "Wonderfull"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 20
$ws.Cells.Item(2, 14).Value = 20000
$ws.Cells.Item(2, 15).Value = 20000
$ws.Cells.Item(2, 16).Value = 20000
$ws.Cells.Item(2, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(2, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(2, 19).Value = 1111
$ws.Cells.Item(2, 20).Value = 18
